# Refresh the nowcast table for 2025Q4: the six most-recent revision dates
# (2025-09-30 .. 2025-12-15) replace the oldest six (2025-06-30 .. 2025-09-15)
# in the top data rows (2-7), with freshly computed revision figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: new "Row" date labels for rows 2-7 (kept as TEXT, not auto-converted dates) ---
$dates = @("2025-09-30", "2025-10-15", "2025-10-30", "2025-11-15", "2025-11-30", "2025-12-15")
for ($i = 0; $i -lt $dates.Length; $i++) {
    $cell = $ws.Cells.Item($i + 2, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $cell.Style = "Normal"
}

# --- Row 2 ---
$ws.Range("B2").Value = [double]"0.2640050393282547"

# --- Row 3 ---
$ws.Range("B3").Value = [double]"0.26592660516450639"
$ws.Range("D3").Value = [double]"-0.0061882170899833604"
$ws.Range("E3").Value = [double]"-0.00022320412839002459"
$ws.Range("F3").Value = [double]"-0.00055552616566332852"
$ws.Range("G3").Value = [double]"6.2282115895190469e-05"
$ws.Range("H3").Value = [double]"3.4785329125562321e-05"
$ws.Range("I3").Value = [double]"6.7029521694531309e-05"
$ws.Range("K3").Value = [double]"3.3598347237151494e-06"

# --- Row 4 ---
$ws.Range("B4").Value = [double]"0.33517607073174582"
$ws.Range("C4").Value = [double]"0.014881312400461101"
$ws.Range("E4").Value = [double]"-2.8237962311511487e-06"
$ws.Range("F4").Value = [double]"2.8686151737493915e-06"
$ws.Range("H4").Value = [double]"0.00015254774874773962"
$ws.Range("I4").Value = [double]"0.00026834011648371563"
$ws.Range("J4").Value = [double]"0.0029184639597784751"
$ws.Range("K4").Value = [double]"-2.3440040024524933e-06"

# --- Row 5 ---
$ws.Range("B5").Value = [double]"0.33304215226794243"
$ws.Range("D5").Value = [double]"-0.0025835712278940602"
$ws.Range("E5").Value = [double]"-0.0010726796200376746"
$ws.Range("F5").Value = [double]"0.0016666786872637303"
$ws.Range("G5").Value = [double]"-0.00042688497973143047"
$ws.Range("H5").Value = [double]"0.00014876736989733285"
$ws.Range("I5").Value = [double]"-1.1233779166790558e-06"
$ws.Range("K5").Value = [double]"0.00032073881531474724"

# --- Row 6 ---
$ws.Range("B6").Value = [double]"0.25137706776486068"
$ws.Range("C6").Value = [double]"-0.0032071702089411028"
$ws.Range("E6").Value = [double]"2.3320783423242207e-05"
$ws.Range("F6").Value = [double]"0.00023768846647498343"
$ws.Range("H6").Value = [double]"0.00015619707267173547"
$ws.Range("I6").Value = [double]"-0.0043025569001637837"
$ws.Range("K6").Value = [double]"-8.2225607600516781e-08"

# --- Row 7 ---
$ws.Range("B7").Value = [double]"0.24776674328504919"
$ws.Range("D7").Value = [double]"-0.0056092261953253142"
$ws.Range("E7").Value = [double]"-0.00053683899659987986"
$ws.Range("F7").Value = [double]"0.0023649710920070127"
$ws.Range("G7").Value = [double]"0.0031954606598681413"
$ws.Range("H7").Value = [double]"0"
$ws.Range("I7").Value = [double]"0"
$ws.Range("K7").Value = [double]"-0.00012957784520367666"

# --- Column widths (auto-fit in the source workbook after the data refresh) ---
$ws.Columns.Item(3).ColumnWidth = 15.25
$ws.Columns.Item(4).ColumnWidth = 15.25
$ws.Columns.Item(7).ColumnWidth = 16.25
$ws.Columns.Item(8).ColumnWidth = 15.65
$ws.Columns.Item(9).ColumnWidth = 15.75
$ws.Columns.Item(10).ColumnWidth = 15.05

Write-Output "nowcast table refreshed for 2025Q4"
